# format socket stream as per configuration
#
# The "seconds from a reference Julian date" note is replaced by a plain
# "Julian date" note, the example timestamp row now stores actual Julian
# date values (instead of plain seconds / 1E16 placeholders) formatted with
# a 2-decimal "0.00;[Red]0.00" number format, and every row gets a new
# "datatype" classification column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (timestamp): switch the alarm/warn bounds from raw seconds to
#     real Julian-date values, and restyle the warn_high/alarm_high pair
#     (G2:H2) with a 2-decimal number format instead of scientific notation.
$ws.Range("D2").Value = 2415020.5
$ws.Range("E2").Value = 2444239.5

$ws.Range("G2:H2").NumberFormat = "0.00;[Red]0.00"
$ws.Range("G2").Value = 2469807.5
$ws.Range("H2").Value = 2524593.5

# The notes cell for the timestamp row no longer spells out the full
# explanation of the units, just the short name.
$ws.Range("J2").Value = "Julian date"

# --- New column K: a "datatype" tag describing every field's value kind.
$ws.Range("K1").Value = "datatype"
$ws.Range("K1").Font.Bold = $true

$ws.Range("K2").Value = "date"
$ws.Range("K3:K10").Value = "number"
$ws.Range("K11").Value = "string"

# Column G is no longer sized to best-fit its content; give it a fixed,
# slightly wider width instead.
$ws.Range("G1").EntireColumn.ColumnWidth = 10.9166666666667

# Move the active selection to where the editor was last working.
$ws.Range("D4").Select()
